# Insert a new weekly price-report row for Ciboulette (Mercado Mayorista Lo
# Valledor de Santiago) just above the existing row 632, shifting the rest
# of the data block (old rows 632-668) down by one to 633-669.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 632; Excel shifts rows 632..668 down to 633..669
# and carries the row-above formatting (date number format on column D).
$ws.Rows.Item(632).Insert()

$ws.Cells.Item(632, 1).Value  = 6
$ws.Cells.Item(632, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(632, 3).Value  = "Metropolitana"
$ws.Cells.Item(632, 4).Value  = 45041
$ws.Cells.Item(632, 5).Value  = 13
$ws.Cells.Item(632, 6).Value  = 100112039
$ws.Cells.Item(632, 7).Value  = "Ciboulette"
$ws.Cells.Item(632, 8).Value  = "Sin especificar"
$ws.Cells.Item(632, 9).Value  = "Primera"
$ws.Cells.Item(632, 10).Value = 930
$ws.Cells.Item(632, 11).Value = 700
$ws.Cells.Item(632, 12).Value = 800
$ws.Cells.Item(632, 13).Value = 753
$ws.Cells.Item(632, 14).Value = "`$/docena de atados"
$ws.Cells.Item(632, 15).Value = "Región Metropolitana"
$ws.Cells.Item(632, 16).Value = 251
$ws.Cells.Item(632, 17).Value = 3
$ws.Cells.Item(632, 18).Value = "Hortaliza"
